# Part1-TodaysStartup.pptx - "adding team gradebook as dp2 makeup"
#
# The real content edit (per the recorded PowerPoint change-tracking entry
# for shape id=4 / creationId {88299B38-4E01-0C8D-B27E-24B46DF57E87} on the
# first slide, sldId 256) is: the placeholder password text box on Slide 1
# had its blank-line placeholder "__________" replaced with the actual
# attendance password "boxandpointers".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(4)

$tr = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(2)
$run = $para.Runs(1)
$run.Text = "boxandpointers"
